$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 73.69369369369369
$ws.Range("A3").Value = 73.33333333333333
$ws.Range("A4").Value = 75.67567567567568
$ws.Range("A5").Value = 69.72972972972973
$ws.Range("A6").Value = 72.61261261261261
$ws.Range("A7").Value = 72.61261261261261
$ws.Range("A8").Value = 80.36036036036036
$ws.Range("A9").Value = 81.08108108108108
$ws.Range("A10").Value = 79.45945945945945
$ws.Range("A11").Value = 80.36036036036036
$ws.Range("A12").Value = 67.56756756756756
$ws.Range("A13").Value = 72.25225225225225
$ws.Range("A14").Value = 76.21621621621621
$ws.Range("A15").Value = 76.3963963963964
$ws.Range("A16").Value = 76.93693693693695
$ws.Range("A17").Value = 64.32432432432432
$ws.Range("A18").Value = 67.56756756756756
$ws.Range("A19").Value = 71.17117117117117
$ws.Range("A20").Value = 77.83783783783784
$ws.Range("A21").Value = 79.0990990990991
$ws.Range("A22").Value = 78.37837837837837
$ws.Range("A23").Value = 66.48648648648648
$ws.Range("A24").Value = 64.32432432432432
$ws.Range("A25").Value = 64.14414414414415
$ws.Range("A26").Value = 75.49549549549549
$ws.Range("A27").Value = 71.89189189189189
$ws.Range("A28").Value = 74.05405405405405
$ws.Range("A29").Value = 72.7927927927928
$ws.Range("A30").Value = 69.54954954954955
$ws.Range("A31").Value = 72.97297297297297
$ws.Range("A32").Value = 75.31531531531532
$ws.Range("A33").Value = 75.85585585585586
$ws.Range("A34").Value = 76.03603603603604
$ws.Range("A35").Value = 69.009009009009
$ws.Range("A36").Value = 69.54954954954955
$ws.Range("A37").Value = 62.88288288288288
$ws.Range("A38").Value = 70.990990990991
$ws.Range("A39").Value = 66.30630630630631
$ws.Range("A40").Value = 66.66666666666666
$ws.Range("A41").Value = 74.41441441441441
$ws.Range("A42").Value = 76.21621621621621
$ws.Range("A43").Value = 76.75675675675676
$ws.Range("A44").Value = 75.49549549549549
$ws.Range("A45").Value = 75.85585585585586
$ws.Range("A46").Value = 75.67567567567568
$ws.Range("A47").Value = 68.46846846846847
$ws.Range("A48").Value = 68.28828828828829
$ws.Range("A49").Value = 73.51351351351352
